# Update "想去人数" (column F) figures across sheets, per gh-pages regeneration.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 7910
$ws.Cells.Item(3, 6).Value = 108
$ws.Cells.Item(5, 6).Value = 15904
$ws.Cells.Item(6, 6).Value = 45
$ws.Cells.Item(9, 6).Value = 453
$ws.Cells.Item(11, 6).Value = 442
$ws.Cells.Item(12, 6).Value = 786
$ws.Cells.Item(15, 6).Value = 349
$ws.Cells.Item(16, 6).Value = 22
$ws.Cells.Item(17, 6).Value = 302
$ws.Cells.Item(18, 6).Value = 142
$ws.Cells.Item(20, 6).Value = 350
$ws.Cells.Item(21, 6).Value = 1102
$ws.Cells.Item(23, 6).Value = 654
$ws.Cells.Item(24, 6).Value = 2242
$ws.Cells.Item(25, 6).Value = 764
$ws.Cells.Item(26, 6).Value = 59
$ws.Cells.Item(27, 6).Value = 563
$ws.Cells.Item(29, 6).Value = 621
$ws.Cells.Item(30, 6).Value = 559

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 67
$ws.Cells.Item(4, 6).Value = 331
$ws.Cells.Item(11, 6).Value = 5

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 477

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 477
$ws.Cells.Item(3, 6).Value = 7910
$ws.Cells.Item(4, 6).Value = 108
$ws.Cells.Item(7, 6).Value = 15904
$ws.Cells.Item(8, 6).Value = 45
$ws.Cells.Item(11, 6).Value = 453
$ws.Cells.Item(12, 6).Value = 67
$ws.Cells.Item(14, 6).Value = 442
$ws.Cells.Item(15, 6).Value = 331
$ws.Cells.Item(18, 6).Value = 786
$ws.Cells.Item(21, 6).Value = 349
$ws.Cells.Item(23, 6).Value = 22
$ws.Cells.Item(27, 6).Value = 302
$ws.Cells.Item(28, 6).Value = 142
$ws.Cells.Item(30, 6).Value = 350
$ws.Cells.Item(31, 6).Value = 1102
$ws.Cells.Item(33, 6).Value = 654
$ws.Cells.Item(34, 6).Value = 2242
$ws.Cells.Item(35, 6).Value = 764
$ws.Cells.Item(36, 6).Value = 59
$ws.Cells.Item(37, 6).Value = 563
$ws.Cells.Item(39, 6).Value = 5
$ws.Cells.Item(40, 6).Value = 621
$ws.Cells.Item(41, 6).Value = 559
